# Organizacion.xlsx update
# "Correcciones leves, actualizacion de ReadMe y Organizacion"
#
# Fills in the remaining task rows (8-12) that were placeholder rows,
# adds 5 extra blank rows below them, widens column B a bit, and moves
# the selection/scroll position to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 12 -> task 8
# ---------------------------------------------------------------------
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Creacion de la estructura static y css base"
$ws.Range("C12").Value = "Sergio"
$ws.Range("D12").Value = "2025-04-04"
$ws.Range("E12").Value = "2025-04-05"
$ws.Range("F12").Value = "✅ Hecho"

# ---------------------------------------------------------------------
# Row 13 -> task 9
# ---------------------------------------------------------------------
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Modificacion de la pagina principal"
$ws.Range("C13").Value = "Samuel"
$ws.Range("D13").Value = "2025-04-05"
$ws.Range("E13").Value = "2025-04-05"
$ws.Range("F13").Value = "✅ Hecho"

# ---------------------------------------------------------------------
# Row 14 -> task 10
# ---------------------------------------------------------------------
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Mejoras en la visualizacion del listado haciendolo mediante tablas añadiendo algunos atributos mas para mejorar la comprension, cambios en el estilo de los botones, cambios en el estilo en los links a los detalles, modificado boton de volver al listado y añadido volver a pagina anterior, modificados los botones de volver de las paginas de creacion edicion y eliminacion "
$ws.Range("C14").Value = "Jon"
$ws.Range("D14").Value = "2025-04-05"
$ws.Range("E14").Value = "2025-04-06"
$ws.Range("F14").Value = "✅ Hecho"

# ---------------------------------------------------------------------
# Row 15 -> task 11
# ---------------------------------------------------------------------
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Modificacion de los botones de eliminacion"
$ws.Range("C15").Value = "Sergio"
$ws.Range("D15").Value = "2025-04-06"
$ws.Range("E15").Value = "2025-04-06"
$ws.Range("F15").Value = "✅ Hecho"

# ---------------------------------------------------------------------
# Row 16 -> task 12 (brand new row)
# ---------------------------------------------------------------------
$ws.Range("A16").Value = 12
$ws.Range("B16").Value = "Correcciones leves, actualizacion ReadMe y Organización"
$ws.Range("C16").Value = "Jon"
$ws.Range("D16").Value = "2025-04-06"
$ws.Range("E16").Value = "2025-04-06"
$ws.Range("F16").Value = "✅ Hecho"

# Row heights for the wrapped-text rows (auto height from the original file)
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 195
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30

# ---------------------------------------------------------------------
# New blank rows 17-21, formatted the same as the data rows above them
# ---------------------------------------------------------------------
$ws.Range("A12:G12").Copy()
$ws.Range("A17:G21").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A17:G21").ClearContents()

# ---------------------------------------------------------------------
# Column B is a little wider now that it holds the Responsable names
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 37

# ---------------------------------------------------------------------
# Move the view / selection to the newly added last row
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B16").Select()
